$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("N2").Value = 3.8
$ws.Range("P2").Value = 1.97
$ws.Range("T2").Value = 1.76
# Row 3
$ws.Range("S3").Value = 3.1
$ws.Range("T3").Value = 1.67
# Row 5
$ws.Range("L5").Value = 1.55
$ws.Range("M5").Value = 1.13
$ws.Range("W5").Value = 1.47
$ws.Range("Y5").Value = 9.800000000000001
$ws.Range("AE5").Value = 65
$ws.Range("AN5").Value = 1000
$ws.Range("AO5").Value = 100
# Row 8
$ws.Range("U8").Value = 1.7
$ws.Range("AB8").Value = 980
$ws.Range("AG8").Value = 980
$ws.Range("AH8").Value = 980
$ws.Range("AI8").Value = 980
# Row 11
$ws.Range("L11").Value = 1.46
$ws.Range("Z11").Value = 980
$ws.Range("AL11").Value = 980
# Row 12
$ws.Range("AE12").Value = 42
$ws.Range("AF12").Value = 16
# Row 13
$ws.Range("T13").Value = 1.6
$ws.Range("U13").Value = 2.36
# Row 15
$ws.Range("G15").Value = 1.96
$ws.Range("R15").Value = 1.43
# Row 16
$ws.Range("F16").Value = 2.74
$ws.Range("G16").Value = 3.1
$ws.Range("K16").Value = 3.75
$ws.Range("M16").Value = 1.07
$ws.Range("N16").Value = 3.45
$ws.Range("P16").Value = 1.79
$ws.Range("Q16").Value = 1.85
$ws.Range("R16").Value = 1.32
$ws.Range("S16").Value = 3.4
$ws.Range("T16").Value = 1.74
$ws.Range("V16").Value = 1.53
$ws.Range("X16").Value = 16.5
$ws.Range("Y16").Value = 13.5
$ws.Range("Z16").Value = 22
$ws.Range("AA16").Value = 50
$ws.Range("AB16").Value = 14
$ws.Range("AC16").Value = 9.6
$ws.Range("AD16").Value = 15
$ws.Range("AE16").Value = 36
$ws.Range("AF16").Value = 24
$ws.Range("AG16").Value = 15.5
$ws.Range("AH16").Value = 21
$ws.Range("AI16").Value = 50
$ws.Range("AJ16").Value = 60
$ws.Range("AK16").Value = 42
$ws.Range("AL16").Value = 55
$ws.Range("AM16").Value = 110
$ws.Range("AN16").Value = 36
$ws.Range("AO16").Value = 30
# Row 17
$ws.Range("G17").Value = 4.5
$ws.Range("J17").Value = 3.4
$ws.Range("L17").Value = 1.26
$ws.Range("N17").Value = 4.7
$ws.Range("R17").Value = 1.5
$ws.Range("S17").Value = 2.38
$ws.Range("T17").Value = 1.59
$ws.Range("U17").Value = 2.34
$ws.Range("V17").Value = 1.85
$ws.Range("W17").Value = 1.31
$ws.Range("X17").Value = 26
$ws.Range("Y17").Value = 15
$ws.Range("Z17").Value = 18
$ws.Range("AA17").Value = 30
$ws.Range("AB17").Value = 22
$ws.Range("AC17").Value = 11.5
$ws.Range("AD17").Value = 13.5
$ws.Range("AE17").Value = 24
$ws.Range("AF17").Value = 36
$ws.Range("AG17").Value = 19
$ws.Range("AH17").Value = 19.5
$ws.Range("AI17").Value = 36
$ws.Range("AJ17").Value = 80
$ws.Range("AK17").Value = 48
$ws.Range("AL17").Value = 50
$ws.Range("AM17").Value = 80
$ws.Range("AN17").Value = 38
$ws.Range("AO17").Value = 13.5
# Row 18
$ws.Range("T18").Value = 1.52
$ws.Range("AB18").Value = 980
$ws.Range("AC18").Value = 980
$ws.Range("AF18").Value = 1000
$ws.Range("AG18").Value = 1000
$ws.Range("AH18").Value = 980
$ws.Range("AJ18").Value = 970
$ws.Range("AK18").Value = 1000
$ws.Range("AL18").Value = 980
# Row 19
$ws.Range("H19").Value = 3.05
$ws.Range("X19").Value = 980
$ws.Range("Y19").Value = 980
$ws.Range("Z19").Value = 980
$ws.Range("AB19").Value = 980
$ws.Range("AC19").Value = 980
$ws.Range("AD19").Value = 980
$ws.Range("AE19").Value = 980
$ws.Range("AF19").Value = 980
$ws.Range("AG19").Value = 980
$ws.Range("AH19").Value = 980
$ws.Range("AI19").Value = 980
$ws.Range("AJ19").Value = 980
$ws.Range("AK19").Value = 980
$ws.Range("AL19").Value = 980
$ws.Range("AN19").Value = 980
$ws.Range("AO19").Value = 980
# Row 20
$ws.Range("Q20").Value = 1.36
$ws.Range("T20").Value = 1.6
$ws.Range("U20").Value = 2
# Row 21
$ws.Range("S21").Value = 2.3
$ws.Range("AI21").Value = 36
$ws.Range("AM21").Value = 55
$ws.Range("AN21").Value = 13
# Row 22
$ws.Range("F22").Value = 5
$ws.Range("G22").Value = 9.6
$ws.Range("H22").Value = 1.4
$ws.Range("I22").Value = 1.56
$ws.Range("N22").Value = 1.1
$ws.Range("P22").Value = 2.88
$ws.Range("Q22").Value = 1.39
$ws.Range("R22").Value = 1.61
$ws.Range("V22").Value = 2.78
$ws.Range("W22").Value = 1.11
# Row 23
$ws.Range("S23").Value = 2.06
# Row 24
$ws.Range("I24").Value = 2.22
$ws.Range("J24").Value = 3.25
$ws.Range("N24").Value = 2.86
$ws.Range("P24").Value = 1.88
$ws.Range("Q24").Value = 1.71
$ws.Range("R24").Value = 1.24
$ws.Range("V24").Value = 1.81
# Row 25
$ws.Range("J25").Value = 5.1
$ws.Range("L25").Value = 1.26
$ws.Range("P25").Value = 2.68
$ws.Range("S25").Value = 2.4
$ws.Range("V25").Value = 1.17
$ws.Range("W25").Value = 2.86
$ws.Range("AA25").Value = 180
$ws.Range("AO25").Value = 70
# Row 26
$ws.Range("Q26").Value = 2.24
$ws.Range("X26").Value = 13.5
$ws.Range("Y26").Value = 11.5
$ws.Range("Z26").Value = 17.5
$ws.Range("AA26").Value = 980
$ws.Range("AB26").Value = 10.5
$ws.Range("AC26").Value = 7.4
$ws.Range("AD26").Value = 13
$ws.Range("AE26").Value = 36
$ws.Range("AF26").Value = 980
$ws.Range("AG26").Value = 14
$ws.Range("AH26").Value = 980
$ws.Range("AI26").Value = 55
$ws.Range("AJ26").Value = 60
$ws.Range("AK26").Value = 980
$ws.Range("AM26").Value = 150
$ws.Range("AN26").Value = 50
$ws.Range("AO26").Value = 38
# Row 27
$ws.Range("G27").Value = 4.6
$ws.Range("I27").Value = 2.74
$ws.Range("W27").Value = 1.28
# Row 28
$ws.Range("N28").Value = 1.1
# Row 29
$ws.Range("G29").Value = 2.92
$ws.Range("L29").Value = 1.5
$ws.Range("Z29").Value = 980
$ws.Range("AB29").Value = 9.4
$ws.Range("AJ29").Value = 980
$ws.Range("AK29").Value = 980
$ws.Range("AN29").Value = 980
$ws.Range("AO29").Value = 980
# Row 30
$ws.Range("AI30").Value = 980
# Row 31
$ws.Range("S31").Value = 4.9
# Row 32
$ws.Range("G32").Value = 6
$ws.Range("J32").Value = 4
$ws.Range("P32").Value = 2.08
$ws.Range("T32").Value = 1.8
$ws.Range("AF32").Value = 980
